# Insert a new price observation row at row 75 (pushes existing rows
# 75-183 down to 76-184, so dimension grows from A1:T183 to A1:T184).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(75).Insert()

$ws.Cells.Item(75, 1).Value  = 10
$ws.Cells.Item(75, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(75, 3).Value  = "La Araucanía"
$ws.Cells.Item(75, 4).Value  = 44580
$ws.Cells.Item(75, 5).Value  = 9
$ws.Cells.Item(75, 6).Value  = "Fruta"
$ws.Cells.Item(75, 7).Value  = 100102
$ws.Cells.Item(75, 8).Value  = "Cítricos"
$ws.Cells.Item(75, 9).Value  = 100102006
$ws.Cells.Item(75, 10).Value = "Pomelo"
$ws.Cells.Item(75, 11).Value = "Start Ruby"
$ws.Cells.Item(75, 12).Value = "Primera"
$ws.Cells.Item(75, 13).Value = 125
$ws.Cells.Item(75, 14).Value = 15000
$ws.Cells.Item(75, 15).Value = 15000
$ws.Cells.Item(75, 16).Value = 15000
$ws.Cells.Item(75, 17).Value = "`$/caja 14 kilos granel"
$ws.Cells.Item(75, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(75, 19).Value = 1071
$ws.Cells.Item(75, 20).Value = 14
